$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 01:31"

# Fix swapped order of Montserrat / Islas Malvinas in the shared strings table
# (row 213 was "Islas Malvinas", row 214 was "Montserrat" -- now swapped)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# Updated Covid-19 stats (paises.xlsx refresh)
$ws.Range("B4").Value = 5610123
$ws.Range("C4").Value = 38708
$ws.Range("D4").Value = 2969244
$ws.Range("E4").Value = 2467223
$ws.Range("G4").Value = 529
$ws.Range("H4").Value = 173656
$ws.Range("B5").Value = 3363235
$ws.Range("C5").Value = 23038
$ws.Range("D5").Value = 2478494
$ws.Range("E5").Value = 776087
$ws.Range("G5").Value = 775
$ws.Range("H5").Value = 108654
$ws.Range("B11").Value = 476660
$ws.Range("C11").Value = 8328
$ws.Range("D11").Value = 301525
$ws.Range("E11").Value = 159763
$ws.Range("G11").Value = 275
$ws.Range("H11").Value = 15372
$ws.Range("B13").Value = 382142
$ws.Range("C13").Value = 1833
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 28646
$ws.Range("B17").Value = 299126
$ws.Range("C17").Value = 4557
$ws.Range("E17").Value = 75462
$ws.Range("G17").Value = 111
$ws.Range("H17").Value = 5814
$ws.Range("B27").Value = 122636
$ws.Range("C27").Value = 549
$ws.Range("D27").Value = 108939
$ws.Range("E27").Value = 4667
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 9030
$ws.Range("B32").Value = 96590
$ws.Range("C32").Value = 115
$ws.Range("D32").Value = 60651
$ws.Range("E32").Value = 30766
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 5173
$ws.Range("B35").Value = 86737
$ws.Range("C35").Value = 428
$ws.Range("D35").Value = 54108
$ws.Range("E35").Value = 31148
$ws.Range("G35").Value = 28
$ws.Range("H35").Value = 1481
$ws.Range("B39").Value = 82543
$ws.Range("C39").Value = 603
$ws.Range("D39").Value = 55845
$ws.Range("E39").Value = 24910
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 1788
$ws.Range("B49").Value = 55667
$ws.Range("C49").Value = 953
$ws.Range("D49").Value = 41196
$ws.Range("E49").Value = 13372
$ws.Range("G49").Value = 11
$ws.Range("H49").Value = 1099
$ws.Range("B52").Value = 49485
$ws.Range("C52").Value = 417
$ws.Range("D52").Value = 36834
$ws.Range("E52").Value = 11674
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 977
$ws.Range("B74").Value = 20202
$ws.Range("C74").Value = 190
$ws.Range("D74").Value = 14622
$ws.Range("E74").Value = 5181
$ws.Range("B86").Value = 10060
$ws.Range("C86").Value = 55
$ws.Range("E86").Value = 942
$ws.Range("B102").Value = 6762
$ws.Range("C102").Value = 61
$ws.Range("D102").Value = 6018
$ws.Range("E102").Value = 587
$ws.Range("B114").Value = 4085
$ws.Range("C114").Value = 50
$ws.Range("D114").Value = 2986
$ws.Range("E114").Value = 1019
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 80
$ws.Range("B122").Value = 3077
$ws.Range("C122").Value = 61
$ws.Range("D122").Value = 2138
$ws.Range("E122").Value = 891
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 48
$ws.Range("B143").Value = 1457
$ws.Range("C143").Value = 17
$ws.Range("D143").Value = 1205
$ws.Range("E143").Value = 212
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 40
$ws.Range("B158").Value = 983
$ws.Range("C158").Value = 21
$ws.Range("E158").Value = 492
$ws.Range("B166").Value = 588
$ws.Range("C166").Value = 36
$ws.Range("E166").Value = 436

# Row 213/214 data follows the swapped country labels above
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

